$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "蓝色光标"
$ws.Range("B2").Value = "蓝色光标"
$ws.Range("C2").Value = "三维通信"

$ws.Range("A3").Value = "航天发展"
$ws.Range("B3").Value = "利欧股份"
$ws.Range("C3").Value = "航天发展"

$ws.Range("A4").Value = "利欧股份"
$ws.Range("B4").Value = "岩山科技"
$ws.Range("C4").Value = "蓝色光标"

$ws.Range("A5").Value = "金风科技"
$ws.Range("B5").Value = "航天发展"
$ws.Range("C5").Value = "利欧股份"

$ws.Range("A6").Value = "中国卫星"
$ws.Range("B6").Value = "金风科技"
$ws.Range("C6").Value = "岩山科技"

$ws.Range("A7").Value = "航天电子"
$ws.Range("B7").Value = "航天电子"
$ws.Range("C7").Value = "金风科技"

$ws.Range("A8").Value = "岩山科技"
$ws.Range("B8").Value = "中国卫星"
$ws.Range("C8").Value = "雷科防务"

$ws.Range("A9").Value = "中国卫通"
$ws.Range("B9").Value = "雷科防务"
$ws.Range("C9").Value = "山子高科"

$ws.Range("A10").Value = "易点天下"
$ws.Range("B10").Value = "中国卫通"
$ws.Range("C10").Value = "海格通信"

$ws.Range("A11").Value = "雷科防务"
$ws.Range("B11").Value = "山子高科"
$ws.Range("C11").Value = "中国卫通"

$ws.Range("A12").Value = "山子高科"
$ws.Range("B12").Value = "雪人集团"
$ws.Range("C12").Value = "雪人集团"

$ws.Range("A13").Value = "三维通信"
$ws.Range("B13").Value = "易点天下"
$ws.Range("C13").Value = "中国卫星"

$ws.Range("A14").Value = "通宇通讯"
$ws.Range("B14").Value = "东方财富"
$ws.Range("C14").Value = "航天电子"

$ws.Range("A15").Value = "雪人集团"
$ws.Range("B15").Value = "三六零"
$ws.Range("C15").Value = "万向钱潮"

$ws.Range("A16").Value = "海格通信"
$ws.Range("B16").Value = "昆仑万维"
$ws.Range("C16").Value = "平潭发展"

$ws.Range("A17").Value = "昆仑万维"
$ws.Range("B17").Value = "三维通信"
$ws.Range("C17").Value = "再升科技"

$ws.Range("A18").Value = "华胜天成"
$ws.Range("B18").Value = "海格通信"
$ws.Range("C18").Value = "乾照光电"

$ws.Range("A19").Value = "乾照光电"
$ws.Range("B19").Value = "领益智造"
$ws.Range("C19").Value = "华胜天成"

$ws.Range("A20").Value = "领益智造"
$ws.Range("B20").Value = "中文在线"
$ws.Range("C20").Value = "领益智造"

$ws.Range("A21").Value = "志特新材"
$ws.Range("B21").Value = "美年健康"
$ws.Range("C21").Value = "银河电子"
